# Daily update at 8 AM UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the "last row" date number format currently on A33 before it is
# overwritten, so it can be moved to the new last row.
$lastRowFormat = $ws.Range("A33").NumberFormat

# Append the new day's data as the new last row.
$ws.Range("A34").Value = 45618
$ws.Range("B34").Value = 86
$ws.Range("C34").Value = 69
$ws.Range("D34").Value = 82

# The new last row (34) gets the special "last row" date format that A33
# used to carry.
$ws.Range("A34").NumberFormat = $lastRowFormat

# A33 is no longer the last row, so give it the same date format used by
# the preceding rows (A2:A32).
$ws.Range("A33").NumberFormat = $ws.Range("A32").NumberFormat
